# Auto-generated edit script applying numeric corrections to H:N profit columns
# across multiple sheets, per the commit diff (scheduled runner sheet update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M11").Value = -333
$ws.Range("K11").Value = 473
$ws.Range("H11").Value = 473
$ws.Range("I11").Value = 473
$ws.Range("I19").Value = 4000
$ws.Range("N19").Value = -5013
$ws.Range("L19").Value = 4663
$ws.Range("J19").Value = 4663
$ws.Range("M19").Value = -3825
$ws.Range("K19").Value = 4000
$ws.Range("H19").Value = 4497.25
$ws.Range("M51").Value = -3514.5
$ws.Range("K51").Value = 3998.5
$ws.Range("I51").Value = 3998.5
$ws.Range("H51").Value = 5669.6
$ws.Range("N103").Value = -3120.5
$ws.Range("L103").Value = 1948.5
$ws.Range("J103").Value = 649.5
$ws.Range("H103").Value = 649.5
$ws.Range("K106").Value = 2404.7693
$ws.Range("I106").Value = 2404.7693
$ws.Range("H106").Value = 2390.1428
$ws.Range("M106").Value = -1773.7693
$ws.Range("K132").Value = 4616.4375
$ws.Range("H132").Value = 5841.1665
$ws.Range("I132").Value = 1538.8125
$ws.Range("N132").Value = -48397.625
$ws.Range("L132").Value = 43337.625
$ws.Range("J132").Value = 14445.875
$ws.Range("M132").Value = -2086.4375
$ws.Range("M137").Value = -966.1361999999999
$ws.Range("K137").Value = 3516.1362
$ws.Range("H137").Value = 1327.2222
$ws.Range("I137").Value = 1172.0454
$ws.Range("N137").Value = -11130
$ws.Range("L137").Value = 6030
$ws.Range("J137").Value = 2010

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M45").Value = -334952.5
$ws.Range("K45").Value = 335329.5
$ws.Range("I45").Value = 335329.5
$ws.Range("H45").Value = 185997.81
$ws.Range("M110").Value = -779.9564999999998
$ws.Range("K110").Value = 2824.9565
$ws.Range("H110").Value = 4270.107
$ws.Range("I110").Value = 2824.9565
$ws.Range("N110").Value = -15007.8
$ws.Range("L110").Value = 10917.8
$ws.Range("J110").Value = 10917.8
$ws.Range("K132").Value = 8848.378499999999
$ws.Range("H132").Value = 4082.3809
$ws.Range("I132").Value = 2949.4595
$ws.Range("M132").Value = -6318.378499999999
$ws.Range("H138").Value = 21000
$ws.Range("N138").Value = -31280
$ws.Range("L138").Value = 21000
$ws.Range("J138").Value = 21000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M20").Value = -2736.9473
$ws.Range("K20").Value = 2983.9473
$ws.Range("H20").Value = 3613.3333
$ws.Range("I20").Value = 2983.9473
$ws.Range("H35").Value = 43744.832
$ws.Range("N35").Value = -44364.832
$ws.Range("J35").Value = 43744.832
$ws.Range("L35").Value = 43744.832
$ws.Range("I94").Value = 1820.9375
$ws.Range("H94").Value = 2236.3044
$ws.Range("M94").Value = -1369.9375
$ws.Range("K94").Value = 1820.9375
$ws.Range("K105").Value = 1133.3334
$ws.Range("H105").Value = 1077.8
$ws.Range("I105").Value = 1133.3334
$ws.Range("M105").Value = 613.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N16").Value = -5736.6
$ws.Range("L16").Value = 5162.6
$ws.Range("J16").Value = 5162.6
$ws.Range("M16").Value = -1857.3333
$ws.Range("K16").Value = 2144.3333
$ws.Range("I16").Value = 2144.3333
$ws.Range("H16").Value = 3732.8948
$ws.Range("I22").Value = 300
$ws.Range("H22").Value = 354.42856
$ws.Range("N22").Value = -1127
$ws.Range("L22").Value = 427
$ws.Range("J22").Value = 427
$ws.Range("M22").Value = 50
$ws.Range("K22").Value = 300
$ws.Range("J31").Value = 7301.625
$ws.Range("M31").Value = -4629.625
$ws.Range("K31").Value = 4924.625
$ws.Range("H31").Value = 6113.125
$ws.Range("I31").Value = 4924.625
$ws.Range("N31").Value = -7891.625
$ws.Range("L31").Value = 7301.625
$ws.Range("K34").Value = 4924.625
$ws.Range("H34").Value = 6113.125
$ws.Range("N34").Value = -7705.625
$ws.Range("I34").Value = 4924.625
$ws.Range("L34").Value = 7301.625
$ws.Range("J34").Value = 7301.625
$ws.Range("M34").Value = -4722.625
$ws.Range("M58").Value = -3245.077
$ws.Range("K58").Value = 3448.077
$ws.Range("I58").Value = 3448.077
$ws.Range("H58").Value = 5241.4
$ws.Range("K102").Value = 30000
$ws.Range("I102").Value = 30000
$ws.Range("H102").Value = 30000
$ws.Range("M102").Value = -27566
$ws.Range("N103").Value = -77132
$ws.Range("L103").Value = 74788
$ws.Range("J103").Value = 74788
$ws.Range("M103").Value = -5659
$ws.Range("K103").Value = 6831
$ws.Range("I103").Value = 6831
$ws.Range("H103").Value = 29483.334
$ws.Range("H105").Value = 2743.7778
$ws.Range("N105").Value = -6827
$ws.Range("L105").Value = 3333
$ws.Range("J105").Value = 3333
$ws.Range("I113").Value = 2144.3333
$ws.Range("N113").Value = -9502.6
$ws.Range("L113").Value = 5162.6
$ws.Range("J113").Value = 5162.6
$ws.Range("M113").Value = 25.66670000000022
$ws.Range("K113").Value = 2144.3333
$ws.Range("H113").Value = 3732.8948
$ws.Range("J122").Value = 5999.6665
$ws.Range("M122").Value = -12094.9
$ws.Range("K122").Value = 14544.9
$ws.Range("H122").Value = 5114
$ws.Range("I122").Value = 4848.3
$ws.Range("N122").Value = -22898.9995
$ws.Range("L122").Value = 17998.9995
$ws.Range("I134").Value = 7041.5835
$ws.Range("M134").Value = -18589.7505
$ws.Range("K134").Value = 21124.7505
$ws.Range("H134").Value = 8299.933999999999
$ws.Range("M136").Value = -7794.231
$ws.Range("K136").Value = 10344.231
$ws.Range("H136").Value = 5241.4
$ws.Range("I136").Value = 3448.077

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K5").Value = 5391.2307
$ws.Range("H5").Value = 1761.7727
$ws.Range("N5").Value = -5356.3334
$ws.Range("I5").Value = 1797.0769
$ws.Range("J5").Value = 1710.7778
$ws.Range("M5").Value = -5279.2307
$ws.Range("L5").Value = 5132.3334
$ws.Range("I26").Value = 400049
$ws.Range("H26").Value = 202164.4
$ws.Range("M26").Value = -1199859
$ws.Range("K26").Value = 1200147
$ws.Range("H38").Value = 470.55554
$ws.Range("N38").Value = -2281.75
$ws.Range("L38").Value = 1587.75
$ws.Range("J38").Value = 529.25
$ws.Range("N113").Value = -6992
$ws.Range("L113").Value = 2652
$ws.Range("J113").Value = 884
$ws.Range("H113").Value = 858.4167
$ws.Range("M126").Value = -31059.001
$ws.Range("K126").Value = 35999.001
$ws.Range("I126").Value = 11999.667
$ws.Range("H126").Value = 12749.75
$ws.Range("K132").Value = 14859
$ws.Range("H132").Value = 1908.8334
$ws.Range("I132").Value = 1651
$ws.Range("M132").Value = -12329
$ws.Range("K135").Value = 16173.6921
$ws.Range("H135").Value = 1761.7727
$ws.Range("I135").Value = 1797.0769
$ws.Range("N135").Value = -20467.0002
$ws.Range("L135").Value = 15397.0002
$ws.Range("J135").Value = 1710.7778
$ws.Range("M135").Value = -13638.6921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N126").Value = -24439.25
$ws.Range("L126").Value = 19499.25
$ws.Range("J126").Value = 6499.75
$ws.Range("H126").Value = 2505187.5
$ws.Range("J131").Value = 49999.5
$ws.Range("H131").Value = 49999.5
$ws.Range("N131").Value = -60079.5
$ws.Range("L131").Value = 49999.5
$ws.Range("K132").Value = 4821
$ws.Range("H132").Value = 3576.4285
$ws.Range("I132").Value = 1607
$ws.Range("M132").Value = -2291

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 1875
$ws.Range("H22").Value = 2114.2856
$ws.Range("M22").Value = -1580
$ws.Range("K22").Value = 1875
$ws.Range("I27").Value = 1875
$ws.Range("M27").Value = -1768
$ws.Range("K27").Value = 1875
$ws.Range("H27").Value = 2114.2856
$ws.Range("H40").Value = 2894.2856
$ws.Range("N40").Value = -3771.5
$ws.Range("I40").Value = 2793.4167
$ws.Range("L40").Value = 3499.5
$ws.Range("J40").Value = 3499.5
$ws.Range("M40").Value = -2657.4167
$ws.Range("K40").Value = 2793.4167
$ws.Range("K61").Value = 113089.11
$ws.Range("H61").Value = 85317.836
$ws.Range("I61").Value = 113089.11
$ws.Range("M61").Value = -112887.11
$ws.Range("I113").Value = 113089.11
$ws.Range("M113").Value = -110919.11
$ws.Range("K113").Value = 113089.11
$ws.Range("H113").Value = 85317.836

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L15").Value = 10000
$ws.Range("J15").Value = 10000
$ws.Range("H15").Value = 9750
$ws.Range("N15").Value = -10576
$ws.Range("L54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("H54").Value = 21170
$ws.Range("N54").Value = $null
$ws.Range("M122").Value = -11695.5289
$ws.Range("K122").Value = 14145.5289
$ws.Range("H122").Value = 4858.4585
$ws.Range("I122").Value = 4715.1763
$ws.Range("K132").Value = 7414.263300000001
$ws.Range("H132").Value = 2634.617
$ws.Range("I132").Value = 2471.4211
$ws.Range("N132").Value = -15031.0001
$ws.Range("L132").Value = 9971.000100000001
$ws.Range("J132").Value = 3323.6667
$ws.Range("M132").Value = -4884.263300000001
